$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition) - update "想去人数" (want-to-go count) column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 320
$ws1.Range("F4").Value = 8360
$ws1.Range("F5").Value = 6091
$ws1.Range("F6").Value = 523
$ws1.Range("F7").Value = 108
$ws1.Range("F8").Value = 20
$ws1.Range("F10").Value = 318
$ws1.Range("F11").Value = 1054

# Sheet "全部类型" (all types) - same updates, mirrored at different rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 320
$ws4.Range("F4").Value = 8360
$ws4.Range("F5").Value = 6091
$ws4.Range("F6").Value = 523
$ws4.Range("F7").Value = 108
$ws4.Range("F8").Value = 20
$ws4.Range("F10").Value = 318
$ws4.Range("F15").Value = 1054
